$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update header row: User_Input/Bot_Response -> Question/Response
$ws.Range("A1").Value = "Question"
$ws.Range("B1").Value = "Response"

# Swap rows 7 and 8 (default/... <-> what is sparrow/appraisal report)
$ws.Range("A7").Value = "what is sparrow"
$ws.Range("B7").Value = "appraisal report"
$ws.Range("A8").Value = "default"
$ws.Range("B8").Value = "I'm sorry, I don't understand that. Can you rephrase?"

# Update the selected cell to match the saved selection state
$ws.Range("B12").Select()
